$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1) with the same formatting as the existing
# header row (bold/centered/bordered style used by H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data values for rows 2-4, columns I and J (plain, unstyled cells
# like the rest of the data rows).
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8
